# "letras pablo con etiqueta"
# Update the AR column (header "valor") for every data row (rows 2-51)
# from 25 to 19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AR2:AR51").Value = 19

# Restore the view/selection state captured in the workbook (cosmetic,
# matches the sheetView selection recorded after the edit).
$ws.Range("AU47").Select()
